$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A83").Value = 0
$ws.Range("A95").Value = 1
$ws.Range("A98").Value = 0
$ws.Range("A196").Value = 1
$ws.Range("A202").Value = 1
$ws.Range("A244").Value = 0
$ws.Range("A249").Value = 1
$ws.Range("A292").Value = 0
$ws.Range("A325").Value = 1
$ws.Range("A333").Value = 0
$ws.Range("A344").Value = 1
$ws.Range("A382").Value = 1
$ws.Range("A454").Value = 1
$ws.Range("A484").Value = 0
$ws.Range("A502:A503").Value = 0
$ws.Range("A505:A515").Value = 0
$ws.Range("A517").Value = 0
$ws.Range("A519:A521").Value = 0
$ws.Range("A523:A528").Value = 0
$ws.Range("A530:A540").Value = 0
$ws.Range("A542").Value = 1
$ws.Range("A545").Value = 0
$ws.Range("A546").Value = 1
$ws.Range("A548").Value = 0
$ws.Range("A551").Value = 1
$ws.Range("A554").Value = 1
$ws.Range("A555:A556").Value = 0
$ws.Range("A561").Value = 0
$ws.Range("A562:A563").Value = 1
$ws.Range("A566").Value = 1
$ws.Range("A569").Value = 0
$ws.Range("A572:A573").Value = 0
$ws.Range("A578").Value = 1
$ws.Range("A580").Value = 1
$ws.Range("A588").Value = 0
$ws.Range("A590").Value = 0
$ws.Range("A601").Value = 0
$ws.Range("A606").Value = 0
$ws.Range("A608").Value = 1
$ws.Range("A617").Value = 1
$ws.Range("A618").Value = 0
$ws.Range("A620").Value = 1
$ws.Range("A622").Value = 0
$ws.Range("A630").Value = 1
$ws.Range("A632").Value = 0
$ws.Range("A639").Value = 0
$ws.Range("A642").Value = 1
$ws.Range("A645").Value = 1
$ws.Range("A651:A652").Value = 0
$ws.Range("A654").Value = 0
$ws.Range("A656").Value = 0
$ws.Range("A657").Value = 1
$ws.Range("A666").Value = 0
$ws.Range("A668").Value = 1
$ws.Range("A672:A673").Value = 0
$ws.Range("A682").Value = 1
$ws.Range("A684:A695").Value = 1
$ws.Range("A698:A702").Value = 1
$ws.Range("A704:A714").Value = 1
$ws.Range("A733").Value = 1
$ws.Range("A868").Value = 1
$ws.Range("A876").Value = 0
$ws.Range("A897").Value = 0
$ws.Range("A902").Value = 0
$ws.Range("A905").Value = 0
$ws.Range("A917").Value = 1
$ws.Range("A925").Value = 0
$ws.Range("A927").Value = 0
$ws.Range("A934").Value = 0
$ws.Range("A950").Value = 1
$ws.Range("A964").Value = 1
$ws.Range("A979").Value = 1
$ws.Range("A982").Value = 1
$ws.Range("A985").Value = 0
$ws.Range("A1020").Value = 0
$ws.Range("A1025").Value = 0
$ws.Range("A1046").Value = 0
$ws.Range("A1054").Value = 1
$ws.Range("A1057:A1058").Value = 1
$ws.Range("A1306").Value = 0
$ws.Range("A1325").Value = 1
$ws.Range("A1337").Value = 0
$ws.Range("A1380").Value = 0
$ws.Range("A1397").Value = 0
$ws.Range("A1413").Value = 0
$ws.Range("A1416").Value = 0
$ws.Range("A1426").Value = 1
$ws.Range("A1556").Value = 0
$ws.Range("A1565").Value = 1
$ws.Range("A1579").Value = 0
$ws.Range("A1608").Value = 0
$ws.Range("A1630").Value = 1
$ws.Range("A1654").Value = 1
$ws.Range("A1695").Value = 0
$ws.Range("A1730").Value = 0
